# This script reproduces the commit "Fruta / hortaliza, semanal":
# three new weekly price records are inserted at rows 570-572 of the
# "Repollo" sheet (pushing the existing rows 570-663 down to 573-666).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 570, shifting all data below
# (old rows 570:663) down to 573:666.
$ws.Rows("570:572").Insert()

# --- New row 570 --------------------------------------------------------
$ws.Range("A570").Value = 11
$ws.Range("B570").Value = "Vega Monumental Concepción"
$ws.Range("C570").Value = "Bíobío"
$ws.Range("D570").Value = 45218
$ws.Range("E570").Value = 8
$ws.Range("F570").Value = 100112006
$ws.Range("G570").Value = "Repollo"
$ws.Range("H570").Value = "Copenhague"
$ws.Range("I570").Value = "Primera"
$ws.Range("J570").Value = 1000
$ws.Range("K570").Value = 900
$ws.Range("L570").Value = 1000
$ws.Range("M570").Value = 950
$ws.Range("N570").Value = "`$/unidad"
$ws.Range("O570").Value = "Región Metropolitana"
$ws.Range("P570").Value = 950
$ws.Range("Q570").Value = 1
$ws.Range("R570").Value = "Hortaliza"

# --- New row 571 --------------------------------------------------------
$ws.Range("A571").Value = 11
$ws.Range("B571").Value = "Vega Monumental Concepción"
$ws.Range("C571").Value = "Bíobío"
$ws.Range("D571").Value = 45218
$ws.Range("E571").Value = 8
$ws.Range("F571").Value = 100112006
$ws.Range("G571").Value = "Repollo"
$ws.Range("H571").Value = "Crespo record"
$ws.Range("I571").Value = "Primera"
$ws.Range("J571").Value = 1000
$ws.Range("K571").Value = 800
$ws.Range("L571").Value = 900
$ws.Range("M571").Value = 850
$ws.Range("N571").Value = "`$/unidad"
$ws.Range("O571").Value = "Región Metropolitana"
$ws.Range("P571").Value = 850
$ws.Range("Q571").Value = 1
$ws.Range("R571").Value = "Hortaliza"

# --- New row 572 --------------------------------------------------------
$ws.Range("A572").Value = 11
$ws.Range("B572").Value = "Vega Monumental Concepción"
$ws.Range("C572").Value = "Bíobío"
$ws.Range("D572").Value = 45218
$ws.Range("E572").Value = 8
$ws.Range("F572").Value = 100112006
$ws.Range("G572").Value = "Repollo"
$ws.Range("H572").Value = "Morada(o)"
$ws.Range("I572").Value = "Primera"
$ws.Range("J572").Value = 500
$ws.Range("K572").Value = 1200
$ws.Range("L572").Value = 1200
$ws.Range("M572").Value = 1200
$ws.Range("N572").Value = "`$/unidad"
$ws.Range("O572").Value = "Región Metropolitana"
$ws.Range("P572").Value = 1200
$ws.Range("Q572").Value = 1
$ws.Range("R572").Value = "Hortaliza"
